$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Capture the width of the column to the left (M) before the insert so the
# newly-inserted column N can inherit it, mirroring Excel's native
# "insert column" behaviour of carrying over the left neighbour's formatting.
$leftWidth = $ws.Columns("M").ColumnWidth

# Insert a new blank column before column N ("Late"); this shifts the old
# N/O/P ("Late"/"Date"/"Disbursement") columns right to O/P/Q.
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $leftWidth

# Make "Repayment schedule" the active sheet/tab and restore the cursor to
# the cell that was selected after the edit.
$ws.Select() | Out-Null
$ws.Range("R5").Select() | Out-Null
